$d = $word.ActiveDocument

# --- Hunk 2: FOREIGN KEY (checkedOutBy) references Person (id), -> split into 3 runs with " ON DELETE SET NULL" ---
$r = $d.Content
[void]$r.Find.Execute(") references Person (id),", $false, $false, $false, $false, $false, $true, 1, $false)
$target = $d.Range($r.Start, $r.End)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>) references Person (id)</w:t></w:r><w:r><w:t xml:space="preserve"> ON DELETE SET NULL</w:t></w:r><w:r><w:t>,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.InsertXML($xml2)

# --- Hunk 3: merge "INSERT INTO Book(...)" / "VALUES(...)" paragraphs into a single parameterized paragraph ---
$r = $d.Content
[void]$r.Find.Execute("INSERT INTO Book( authorFName, authorLName, title, genre, description, addedBy, shelf ) ", $false, $false, $false, $false, $false, $true, 1, $false)
$start = $r.Start
$r = $d.Content
[void]$r.Find.Execute("VALUES ('J.K.',  'Rowling',  'Harry Potter',  'Youth',  'Book about wizards', 5, 5)", $false, $false, $false, $false, $false, $true, 1, $false)
$end = $r.End
$target = $d.Range($start, $end)
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002C5085" w:rsidRDefault="002C5085" w:rsidP="002C5085"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>INSERT INTO Book(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>authorFName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>authorLName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, title, genre, description, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>addedBy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, shelf) </w:t></w:r><w:r><w:t>VALUES(?, ?, ?, ?, ?, (SELECT id FROM Librarian WHERE username= ?), ?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.InsertXML($xml3)

# --- Hunk 4: move bookmark, add "Find shelf..." query X, add new "Get number of books checked out" section ---
$r = $d.Content
[void]$r.Find.Execute("Find shelf to put book on", $false, $false, $false, $false, $false, $true, 1, $false)
$start = $r.Start
$end = $d.Content.End
$target = $d.Range($start, $end)
$xml4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="007A049E" w:rsidRPr="007A049E" w:rsidRDefault="007A049E" w:rsidP="00FD61B6"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Find shelf to put book on</w:t></w:r></w:p><w:p w:rsidR="007A049E" w:rsidRPr="00FD61B6" w:rsidRDefault="007A049E" w:rsidP="00FD61B6"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r w:rsidRPr="007A049E"><w:t>SELECT shelf FROM BOOK WHERE genre=</w:t></w:r><w:r><w:t>X</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Get number of books checked out</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>SELECT id FROM Book</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checkedOutBy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramStart"/><w:r><w:t>=(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">SELECT id FROM Person WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pinNum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=</w:t></w:r><w:r><w:t>X)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$target.InsertXML($xml4)
